$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# --- New BOM line note: LCD driver diode is now a throughhole part ---
$note = $ws.Range("K40")
$note.Value = "throughhole: 1N5817-TPCT-ND"
$note.NumberFormat = "@"
$note.HorizontalAlignment = -4108   # xlCenter
$note.VerticalAlignment = -4108     # xlCenter
$note.WrapText = $true

# Row 40 grows to fit the two-line note
$ws.Rows.Item(40).RowHeight = 24

# NOTES column widened to accommodate the new note text
$ws.Columns.Item(11).ColumnWidth = 12.86

# Move the active selection/view to the edited cell
$ws.Range("A22").Select()
$note.Select()
